$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in the Price column (D) that are plain decimal-looking strings get
# auto-converted to numbers by Excel, which loses trailing zeros / dot-grouping
# formatting (e.g. "1.00" -> 1, "0.0810" -> 0.081, "38.280.35" stays text anyway).
# Force text storage by setting NumberFormat to "@" before assigning, then reset
# the cell style back to Normal so no stray style index is left on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "38.280.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.13%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.069.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.67%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.91%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.621"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.23%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.39"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +11.10%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  +4.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0810"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.21%  "

$ws.Range("E11").Value = "  +2.18%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.02"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.83%  "

$ws.Range("E13").Value = "  +2.66%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.46"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.44%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.770"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.43%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.33"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.43%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.080.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.44%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "38.230.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.29%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.52%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.10%  "

$ws.Range("E21").Value = "  +3.49%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.54%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.46"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.71%  "

$ws.Range("E25").Value = "  +3.92%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.59%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.58%  "

$ws.Range("E28").Value = "  +6.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.71%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.32"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.17%  "

$ws.Range("E31").Value = "  +2.68%  "

$ws.Range("E32").Value = "  +4.13%  "

$ws.Range("E33").Value = "  +4.89%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0611"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.82%  "

$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.05"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.94%  "

$ws.Range("E36").Value = "  +0.43%  "

$ws.Range("E37").Value = "  +16.43%  "

$ws.Range("E38").Value = "  +6.89%  "

$ws.Range("E39").Value = "  +0.11%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.530.08"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.58%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.35"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.96%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.79"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.67%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0218"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.07%  "

$ws.Range("E44").Value = "  +4.12%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0928"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.20%  "

$ws.Range("E46").Value = "  +2.00%  "

$ws.Range("E47").Value = "  -1.70%  "

$ws.Range("E48").Value = "  +2.78%  "

$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.15"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.54%  "

$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.95"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.59%  "

$ws.Range("E51").Value = "  +2.85%  "

